$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct typo in site distance matrix: I6 was 11767, should be 1176
# (matches the symmetric value at F9)
$ws.Range("I6").Value = 1176

# Leave selection on the corrected cell, as in the saved workbook
$ws.Range("I7").Select()
